# Handback status report generation: update the "Generate Date" /
# "Handoff Datetime" / "Handback Datetime" timestamps to the latest run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date for the first file.
$wsOverview.Range("G2").Value = "2016-08-29 15:16:49"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
# for the first file.
$wsZhCn.Range("H2").Value = "2016-08-29 15:16:44"
$wsZhCn.Range("K2").Value = "2016-08-29 15:17:07"

# de-de sheet: Correspond Handoff Datetime for the first file.
$wsDeDe.Range("H2").Value = "2016-08-29 15:17:21"
